$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert two new columns at F:G (pushes old F..S to H..U), Excel-native
#    behaviour: inherits number format / border from the column to the left
#    for the cells in the new columns, and shifts col width definitions.
# ---------------------------------------------------------------------------
$ws.Columns("F:G").Insert()

# The insert leaves F2/G2 blank (style copied from column E) - restore the
# text that mirrors the "Facturas pendientes" / "Facturas de la fianza"
# columns (D2/E2) for the new duplicate block.
$ws.Range("F2").Value = $ws.Range("D2").Value2
$ws.Range("G2").Value = $ws.Range("E2").Value2

# ---------------------------------------------------------------------------
# 2) Build the two-column "group caption" header styles used on row 1:
#    - a plain (no border) centerContinuous variant  -> E1 / (target xf #12)
#    - a left-bordered centerContinuous variant       -> D1 & F1 (xf #13)
#    - a right-bordered centerContinuous variant       -> G1 (xf #14)
#    All three share: numFmt 0.0 (164), bold white Century Gothic font,
#    solid teal FF009288 fill - i.e. the same font/numberformat as the
#    existing teal header style, just a different (darker) teal fill.
# ---------------------------------------------------------------------------

# Base the new styles off the existing teal "financial header" cell (C2),
# which already carries fontId=5 / numFmtId=164 - this reuses that font
# instead of re-creating an equivalent one.
$ws.Range("C2").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Interior.Color = 8950272
$ws.Range("E1").HorizontalAlignment = 7

$ws.Range("E1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Borders(7).LineStyle = 1
$ws.Range("D1").Borders(7).Weight = 2

$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("E1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Borders(10).LineStyle = 1
$ws.Range("G1").Borders(10).Weight = 2

# ---------------------------------------------------------------------------
# 3) Captions for the new grouped header cells.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Cotown"
$ws.Range("F1").Value = "Propietario"

$excel.CutCopyMode = 0
